# Update gh-pages to output generated at 456a3b4
# Applies the numeric "want-to-go" count bumps (and one event's time/cover
# image correction) across the 展览 / 演出 / 本地生活 / 全部类型 sheets.

$wb = $excel.ActiveWorkbook

# ---- Sheet 1: 展览 (Exhibitions) ----
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 12836
$ws1.Range("F3").Value = 7180
$ws1.Range("F10").Value = 1008
$ws1.Range("E11").Value = "2024.08.22 11:00-08.24 17:30"
$ws1.Range("I11").Value = "//i0.hdslb.com/bfs/openplatform/202408/bHSLBXrv1723788746248.png"
$ws1.Range("F13").Value = 1017
$ws1.Range("F14").Value = 5
$ws1.Range("F18").Value = 252
$ws1.Range("F19").Value = 371
$ws1.Range("F21").Value = 280
$ws1.Range("F24").Value = 172
$ws1.Range("F25").Value = 375
$ws1.Range("F26").Value = 5243
$ws1.Range("F28").Value = 1434
$ws1.Range("F30").Value = 1382
$ws1.Range("F31").Value = 64
$ws1.Range("F32").Value = 48
$ws1.Range("F33").Value = 1367
$ws1.Range("F36").Value = 598
$ws1.Range("F38").Value = 3739

# ---- Sheet 2: 演出 (Performances) ----
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F4").Value = 3745
$ws2.Range("F5").Value = 3745
$ws2.Range("F14").Value = 11

# ---- Sheet 3: 本地生活 (Local Life) ----
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 9285
$ws3.Range("F3").Value = 561
$ws3.Range("F4").Value = 2025

# ---- Sheet 4: 全部类型 (All Types, combined) ----
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 9285
$ws4.Range("F3").Value = 561
$ws4.Range("F4").Value = 2025
$ws4.Range("F5").Value = 12836
$ws4.Range("F6").Value = 7180
$ws4.Range("F8").Value = 3745
$ws4.Range("F10").Value = 1008
$ws4.Range("E11").Value = "2024.08.22 11:00-08.24 17:30"
$ws4.Range("I11").Value = "//i0.hdslb.com/bfs/openplatform/202408/bHSLBXrv1723788746248.png"
$ws4.Range("F13").Value = 1017
$ws4.Range("F14").Value = 5
$ws4.Range("F18").Value = 252
$ws4.Range("F19").Value = 371
$ws4.Range("F21").Value = 280
$ws4.Range("F27").Value = 172
$ws4.Range("F28").Value = 375
$ws4.Range("F29").Value = 5243
$ws4.Range("F31").Value = 1434
$ws4.Range("F36").Value = 1382
$ws4.Range("F37").Value = 64
$ws4.Range("F38").Value = 1367
$ws4.Range("F40").Value = 598
$ws4.Range("F41").Value = 11
$ws4.Range("F47").Value = 3739
